$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 11) duplicating the data already present in row 10.
$ws.Range("A11").Value = "Mặt sau"
$ws.Range("B11").Value = "PHAM DUY LONG"
$ws.Range("C11").Value = "S Trà Co, Thanh Cái, Qung NInh phó Móng Khu Trang Ginl Trà Co, Thanh Móng Cál, phó"

# D11 and E11 look numeric/date-like ("03/12/2006", "022206004066") but must
# stay plain text (leading zeros / non-date literal), so force text format
# before assigning, then restore the default style so no style index sticks.
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "03/12/2006"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "022206004066"
$ws.Range("E11").Style = "Normal"

$ws.Range("F11").Value = "0v12/2031"
$ws.Range("G11").Value = "Việt Nam"
$ws.Range("H11").Value = "Hải Xuan, Thành phố Móng Cái, Quảng Ninh Hải Xuán, Thành phó Móng Cá"
$ws.Range("I11").Value = "Nam"
